# Insert a new data row at row 54 (pushing existing rows 54..122 down to 55..123)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 54; this shifts rows 54:122 down to 55:123
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with the new record
$ws.Range("A54").Value2 = 11
$ws.Range("B54").Value2 = "Vega Monumental Concepción"
$ws.Range("C54").Value2 = "Bíobío"
$ws.Range("D54").Value2 = 44671
$ws.Range("E54").Value2 = 8
$ws.Range("F54").Value2 = 100112043
$ws.Range("G54").Value2 = "Pepino ensalada"
$ws.Range("H54").Value2 = "Sin especificar"
$ws.Range("I54").Value2 = "Primera"
$ws.Range("J54").Value2 = 180
$ws.Range("K54").Value2 = 13000
$ws.Range("L54").Value2 = 14000
$ws.Range("M54").Value2 = 13444
$ws.Range("N54").Value2 = "`$/caja 60 unidades"
$ws.Range("O54").Value2 = "Región Metropolitana"
$ws.Range("P54").Value2 = 224
$ws.Range("Q54").Value2 = 60
$ws.Range("R54").Value2 = "Hortaliza"

# Ensure the date cell keeps the same date/time number format as the other date cells (style index 2)
$ws.Range("D54").NumberFormat = $ws.Range("D55").NumberFormat
